# Remove the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph
# together with the blank paragraph before it, the blank paragraph
# after it, and the following (blank, page-break-before) paragraph —
# i.e. the four paragraphs that used to sit between the "Requisitos"
# list and the trailing page-break paragraphs.

$d = $word.ActiveDocument

# Locate the anchor text; narrows $anchor to the matched range.
$anchor = $d.Content
$anchor.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false,
                      $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target = $anchor.Start

# Resolve the paragraph index that contains the matched text.
$idx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $target -and $target -lt $p.Range.End) {
        $idx = $i
        break
    }
}

# Delete the blank paragraph before it, the "Ver no Jupiter..." paragraph
# itself, the blank paragraph after it, and the blank page-break paragraph
# that follows — four paragraphs in total — leaving the remaining blank
# paragraph and the final page-break paragraph untouched.
$deleteStart = $d.Paragraphs.Item($idx - 1).Range.Start
$deleteEnd = $d.Paragraphs.Item($idx + 3).Range.Start
$d.Range($deleteStart, $deleteEnd).Delete()
